$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.632.23"
$ws.Range("E2").Value = "  +4.81%  "
$ws.Range("D3").Value = "2.493.17"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.87"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.16"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.30"
$ws.Range("E10").Value = "  +7.37%  "
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.35"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.17"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").Value = "2.878.73"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "2.488.03"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "47.495.51"
$ws.Range("E18").Value = "  +4.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.75"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.34"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.21"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.12"
$ws.Range("E30").Value = "  +6.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  +7.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.39"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.81"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0784"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  +5.67%  "
$ws.Range("E38").Value = "  +4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("E39").Value = "  +4.50%  "
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "122.19"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.39"
$ws.Range("E43").Value = "  +3.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0297"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "1.967.96"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.80"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("E50").Value = "  +11.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.51"
$ws.Range("E51").Value = "  +3.85%  "
